# The new weekly price record is inserted as the first row of this
# Hortaliza/Brócoli block (row 324). Excel shifts every subsequent
# record down by one row, which is exactly what Rows.Item(324).Insert()
# does, and the previously-last record (old row 389) ends up at row 390.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(324).Insert()

$ws.Cells.Item(324, 1).Value = 4
$ws.Cells.Item(324, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(324, 3).Value = 'Los Lagos'
$ws.Cells.Item(324, 4).Value = 44798
$ws.Cells.Item(324, 5).Value = 10
$ws.Cells.Item(324, 6).Value = 100112023
$ws.Cells.Item(324, 7).Value = 'Brócoli'
$ws.Cells.Item(324, 8).Value = 'Sin especificar'
$ws.Cells.Item(324, 9).Value = 'Primera'
$ws.Cells.Item(324, 10).Value = 500
$ws.Cells.Item(324, 11).Value = 1500
$ws.Cells.Item(324, 12).Value = 1500
$ws.Cells.Item(324, 13).Value = 1500
$ws.Cells.Item(324, 14).Value = '$/unidad'
$ws.Cells.Item(324, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(324, 16).Value = 1500
$ws.Cells.Item(324, 17).Value = 1
$ws.Cells.Item(324, 18).Value = 'Hortaliza'
